# Update the NATMI ligand/receptor TPM-derived statistics (Vcam1-Itgb1).
#
# Columns G,H,I,J (ligand average/total expression + derived specificities)
# depend only on the "Sending cluster" (column A).
# Columns M,N,O,P (receptor average/total expression + derived specificities)
# depend only on the "Target cluster" (column D).
# Columns Q,R,S,T (edge weights/specificities) are simply the pairwise
# products: Q = G*M, R = H*N, S = I*O, T = J*P.
#
# This script re-derives all of G,H,I,J,M,N,O,P per row from the
# sending/target cluster and recomputes Q,R,S,T accordingly, matching the
# updated TPM values from the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-cluster ligand stats: avg expr, total expr, specificity(avg), specificity(total)
$ligandByA = @{
    "ECs"              = @(17.93632866666666,  53.808986,          0.1226979812530711,  0.1347750935001359)
    "FAPs"             = @(57.44330666666667,  172.32992,          0.3929554311523962,  0.4316338739568692)
    "Inflammatory-Mac" = @(9.626273333333334,  28.87882,           0.06585095126993876, 0.07233263354328205)
    "MuSCs"            = @(39.29803649999999,  78.59607299999999, 0.2688281328564436,  0.1968591842135532)
    "Resolving-Mac"    = @(21.87880766666666,  65.63642299999999, 0.1496675034681502,  0.1643992147861598)
}

# New per-cluster receptor stats: avg expr, total expr, specificity(avg), specificity(total)
$receptorByD = @{
    "ECs"              = @(121.928739,         365.786217, 0.2282232151508951, 0.2419720431319445)
    "FAPs"             = @(147.91433,          443.74299,  0.2768624053389947, 0.2935413991166814)
    "Inflammatory-Mac" = @(83.50496933333334,  250.514908, 0.1563025480180701, 0.1657186665504434)
    "MuSCs"            = @(91.06846250000001,  182.136925, 0.1704597085236707, 0.1204857969594293)
    "Resolving-Mac"    = @(89.83562999999999,  269.50689,  0.1681521229683693, 0.1782820942415013)
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if (-not $lastRow -or $lastRow -lt 2) { $lastRow = 26 }

for ($r = 2; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()   # Sending cluster
    $d = $ws.Cells.Item($r, 4).Value()   # Target cluster

    if (-not $ligandByA.ContainsKey($a) -or -not $receptorByD.ContainsKey($d)) {
        continue
    }

    $lig = $ligandByA[$a]
    $rec = $receptorByD[$d]

    $g = $lig[0]; $h = $lig[1]; $i = $lig[2]; $j = $lig[3]
    $m = $rec[0]; $n = $rec[1]; $o = $rec[2]; $p = $rec[3]

    $ws.Cells.Item($r, 7).Value  = $g    # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $h    # H: Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $i    # I: Ligand derived specificity (average)
    $ws.Cells.Item($r, 10).Value = $j    # J: Ligand derived specificity (total)

    $ws.Cells.Item($r, 13).Value = $m    # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $n    # N: Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $o    # O: Receptor derived specificity (average)
    $ws.Cells.Item($r, 16).Value = $p    # P: Receptor derived specificity (total)

    $ws.Cells.Item($r, 17).Value = $g * $m   # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $h * $n   # R: Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $i * $o   # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $j * $p   # T: Edge total expression derived specificity
}
